# Auto-generated from the authoritative diff; updates the
# cryptocurrency price/volume snapshot cells in columns B:E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.217.27"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "2.487.07"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'567.27"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").Value = "'165.59"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.80%  "

$ws.Range("E9").Value = "  -0.87%  "

$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("E11").Value = "  -2.95%  "

$ws.Range("D12").Value = "'4.86"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").Value = "2.943.28"
$ws.Range("E13").Value = "  -1.02%  "

$ws.Range("D14").Value = "69.175.69"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").Value = "'24.04"
$ws.Range("E16").Value = "  -3.13%  "

$ws.Range("D17").Value = "2.529.86"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "'11.14"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").Value = "'352.72"
$ws.Range("E19").Value = "  +0.77%  "

$ws.Range("D20").Value = "'7.33"
$ws.Range("E20").Value = "  -3.23%  "

$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("D22").Value = "'1.91"
$ws.Range("E22").Value = "  -3.41%  "

$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("D24").Value = "'69.10"
$ws.Range("E24").Value = "  -1.74%  "

$ws.Range("E25").Value = "  -3.62%  "

$ws.Range("E26").Value = "  -0.77%  "

$ws.Range("D27").Value = "'8.59"
$ws.Range("E27").Value = "  -3.46%  "

$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("E29").Value = "  -2.56%  "

$ws.Range("B30").Value = "POPCAT"
$ws.Range("C30").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Range("D30").Value = "'3.63"
$ws.Range("E30").Value = "  +142.00%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'7.46"
$ws.Range("E31").Value = "  -4.28%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.19"
$ws.Range("E32").Value = "  -3.89%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'436.64"
$ws.Range("E33").Value = "  -5.63%  "

$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("D36").Value = "'154.45"
$ws.Range("E36").Value = "  -1.81%  "

$ws.Range("E37").Value = "  -3.89%  "

$ws.Range("D38").Value = "'19.04"
$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("D39").Value = "'18.07"
$ws.Range("E39").Value = "  -2.16%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").Value = "'0.311"
$ws.Range("E41").Value = "  -2.11%  "

$ws.Range("D42").Value = "'4.56"
$ws.Range("E42").Value = "  -2.82%  "

$ws.Range("E43").Value = "  -2.38%  "

$ws.Range("E44").Value = "  -2.04%  "

$ws.Range("D45").Value = "'1.06"
$ws.Range("E45").Value = "  -4.46%  "

$ws.Range("D46").Value = "'138.09"
$ws.Range("E46").Value = "  -2.67%  "

$ws.Range("D47").Value = "'3.41"
$ws.Range("E47").Value = "  -1.50%  "

$ws.Range("D48").Value = "'0.502"
$ws.Range("E48").Value = "  -3.32%  "

$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").Value = "'0.0926"
$ws.Range("E51").Value = "  -0.40%  "
